$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember column H's width so the newly inserted column can match it
$colHWidth = $ws.Columns(8).ColumnWidth

# Insert a new column before column I (9th column), shifting cX0's
# concentration columns, etc. to the right.
$ws.Columns("I:I").Insert()

# New column I should keep the same width as column H (11 characters)
$ws.Columns(9).ColumnWidth = $colHWidth

# Header for the new column: "mX0" (mass derived from cX0 concentration)
$ws.Cells.Item(1, 9).Value = "mX0"

# Units row for the new column: "g" (grams), matching H2's companion unit
$ws.Cells.Item(2, 9).Value = "g"

# Apply the same numeric style as column H (0.000 format) to the new column's data rows
$ws.Cells.Item(3, 9).Style = $ws.Cells.Item(3, 8).Style
$ws.Cells.Item(4, 9).Style = $ws.Cells.Item(4, 8).Style

# Formulas: I3 = H3*0.5, I4 shares the same formula pattern
$ws.Range("I3:I4").FormulaR1C1 = "=RC[-1]*0.5"

# Update the selection to match the post-edit state
$ws.Range("I3:I4").Select()
